$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.420.78"

# Row 3
$ws.Range("D3").Value = "3.497.97"
$ws.Range("E3").Value = "  +3.84%  "

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.NumberFormat = "General"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "585.67"
$c.NumberFormat = "General"
$ws.Range("E5").Value = "  +2.78%  "

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "147.70"
$c.NumberFormat = "General"
$ws.Range("E6").Value = "  +6.52%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("E8").Value = "  +1.36%  "

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "7.72"
$c.NumberFormat = "General"

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.127"
$c.NumberFormat = "General"
$ws.Range("E10").Value = "  +4.53%  "

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.400"
$c.NumberFormat = "General"
$ws.Range("E11").Value = "  +4.58%  "

# Row 12
$ws.Range("D12").Value = "4.095.53"
$ws.Range("E12").Value = "  +3.75%  "

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "29.93"
$c.NumberFormat = "General"
$ws.Range("E13").Value = "  +8.04%  "

# Row 14
$ws.Range("E14").Value = "  -0.35%  "

# Row 15
$ws.Range("D15").Value = "3.500.49"
$ws.Range("E15").Value = "  +3.72%  "

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0000175"
$c.NumberFormat = "General"
$ws.Range("E16").Value = "  +4.33%  "

# Row 17
$ws.Range("D17").Value = "63.410.59"
$ws.Range("E17").Value = "  +3.96%  "

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "6.28"
$c.NumberFormat = "General"
$ws.Range("E18").Value = "  +3.27%  "

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "14.29"
$c.NumberFormat = "General"

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "9.50"
$c.NumberFormat = "General"
$ws.Range("E20").Value = "  +7.03%  "

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "395.49"
$c.NumberFormat = "General"
$ws.Range("E21").Value = "  +3.60%  "

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.566"
$c.NumberFormat = "General"
$ws.Range("E22").Value = "  +3.22%  "

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "75.47"
$c.NumberFormat = "General"
$ws.Range("E23").Value = "  -0.06%  "

# Row 24
$ws.Range("E24").Value = "  -0.02%  "

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.0000121"
$c.NumberFormat = "General"
$ws.Range("E25").Value = "  +9.31%  "

# Row 26
$ws.Range("D26").Value = "3.644.14"
$ws.Range("E26").Value = "  +3.97%  "

# Row 27
$ws.Range("E27").Value = "  -1.32%  "

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.83"
$c.NumberFormat = "General"
$ws.Range("E28").Value = "  +9.43%  "

# Row 29
$ws.Range("E29").Value = "  -0.40%  "

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "8.26"
$c.NumberFormat = "General"
$ws.Range("E30").Value = "  +5.84%  "

# Row 31
$ws.Range("E31").Value = "  +2.91%  "

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.43"
$c.NumberFormat = "General"
$ws.Range("E32").Value = "  +6.89%  "

# Row 33
$ws.Range("E33").Value = "  +0.00%  "

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "23.87"
$c.NumberFormat = "General"
$ws.Range("E34").Value = "  +4.27%  "

# Row 35
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "7.20"
$c.NumberFormat = "General"
$ws.Range("E35").Value = "  +5.09%  "

# Row 36
$ws.Range("B36").Value = "EnergySwap"
$ws.Range("C36").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "32.58"
$c.NumberFormat = "General"
$ws.Range("E36").Value = "  +28.57%  "

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.35"
$c.NumberFormat = "General"
$ws.Range("E37").Value = "  +8.82%  "

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "172.43"
$c.NumberFormat = "General"
$ws.Range("E38").Value = "  +3.12%  "

# Row 39
$ws.Range("E39").Value = "  +9.33%  "

# Row 40
$ws.Range("D40").Value = "3.532.54"
$ws.Range("E40").Value = "  +3.72%  "

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.0771"
$c.NumberFormat = "General"
$ws.Range("E41").Value = "  +1.56%  "

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.804"
$c.NumberFormat = "General"
$ws.Range("E42").Value = "  +3.99%  "

# Row 43
$ws.Range("E43").Value = "  +7.84%  "

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "4.52"
$c.NumberFormat = "General"

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "42.54"
$c.NumberFormat = "General"
$ws.Range("E45").Value = "  +0.17%  "

# Row 46
$ws.Range("E46").Value = "  +10.37%  "

# Row 47
$ws.Range("D47").Value = "2.615.95"
$ws.Range("E47").Value = "  +6.57%  "

# Row 48
$ws.Range("E48").Value = "  +7.27%  "

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.29"
$c.NumberFormat = "General"
$ws.Range("E49").Value = "  +13.69%  "

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "6.78"
$c.NumberFormat = "General"
$ws.Range("E50").Value = "  +2.64%  "

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0271"
$c.NumberFormat = "General"
$ws.Range("E51").Value = "  +5.13%  "
